$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes what was row 3 (for A, E, F, G, H, Q, R), with a new value for B
$ws.Range("A2").Value = 112181658
$ws.Range("B2").Value = 78699
$ws.Range("E2").Value = 6458
$ws.Range("F2").Value = "Lunglav"
$ws.Range("G2").Value = "Lobaria pulmonaria"
$ws.Range("H2").Value = "(L.) Hoffm."
$ws.Range("Q2").Value = 667994
$ws.Range("R2").Value = 7183150

# Row 3 becomes what was row 2 (for A, E, F, G, H, Q, R), with a new value for B
$ws.Range("A3").Value = 112181823
$ws.Range("B3").Value = 77389
$ws.Range("E3").Value = 228912
$ws.Range("F3").Value = "Mörk kolflarnlav"
$ws.Range("G3").Value = "Carbonicola myrmecina"
$ws.Range("H3").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q3").Value = 667988
$ws.Range("R3").Value = 7183053
